# Scheduled-runner refresh of market/profit figures (currentAveragePrice*,
# LevePrice*, LeveProfit*) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
# Values below are the freshly recomputed figures; a couple of rows also
# flip which of the NQ/HQ profit columns (M/N) is populated, so those cells
# are cleared outright rather than merely zeroed.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 42509.5
$ws.Range("I21").Value = 42509.5
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 42509.5
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -42041.5
$ws.Range("N21").ClearContents()

$ws.Range("H23").Value = 42509.5
$ws.Range("I23").Value = 42509.5
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 42509.5
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -42275.5
$ws.Range("N23").ClearContents()

$ws.Range("H29").Value = 1833
$ws.Range("I29").Value = 299.5
$ws.Range("J29").Value = 4900
$ws.Range("K29").Value = 898.5
$ws.Range("L29").Value = 14700
$ws.Range("M29").Value = -617.5
$ws.Range("N29").Value = -15262

$ws.Range("H38").Value = 437.06668
$ws.Range("I38").Value = 254
$ws.Range("J38").Value = 3000
$ws.Range("K38").Value = 762
$ws.Range("L38").Value = 9000
$ws.Range("M38").Value = -390
$ws.Range("N38").Value = -9744

$ws.Range("H58").Value = 72921.78999999999
$ws.Range("I58").Value = 497.5
$ws.Range("J58").Value = 127240
$ws.Range("K58").Value = 1492.5
$ws.Range("L58").Value = 381720
$ws.Range("M58").Value = -1342.5
$ws.Range("N58").Value = -382020

$ws.Range("H62").Value = 4021.6667
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 4021.6667
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 4021.6667
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -5269.6667

$ws.Range("H65").Value = 4021.6667
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 4021.6667
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 20108.3335
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -26348.3335

$ws.Range("H87").Value = 28474.875
$ws.Range("J87").Value = 28474.875
$ws.Range("L87").Value = 28474.875
$ws.Range("N87").Value = -30970.875

$ws.Range("H90").Value = 28474.875
$ws.Range("J90").Value = 28474.875
$ws.Range("L90").Value = 85424.625
$ws.Range("N90").Value = -97904.625

$ws.Range("H141").Value = 492477.06
$ws.Range("I141").Value = 1743
$ws.Range("J141").Value = 983211.1
$ws.Range("K141").Value = 5229
$ws.Range("L141").Value = 2949633.3
$ws.Range("M141").Value = -49
$ws.Range("N141").Value = -2959993.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1862.75
$ws.Range("I45").Value = 1100
$ws.Range("J45").Value = 2625.5
$ws.Range("K45").Value = 1100
$ws.Range("L45").Value = 2625.5
$ws.Range("M45").Value = -723
$ws.Range("N45").Value = -3379.5

$ws.Range("H63").Value = 2619.3
$ws.Range("I63").Value = 1825.9333
$ws.Range("K63").Value = 1825.9333
$ws.Range("M63").Value = -1139.9333

$ws.Range("H66").Value = 2619.3
$ws.Range("I66").Value = 1825.9333
$ws.Range("K66").Value = 9129.666499999999
$ws.Range("M66").Value = -5697.666499999999

$ws.Range("H74").Value = 1743.1666
$ws.Range("I74").Value = 2029.1428
$ws.Range("K74").Value = 2029.1428
$ws.Range("M74").Value = -1155.1428

$ws.Range("H77").Value = 1743.1666
$ws.Range("I77").Value = 2029.1428
$ws.Range("K77").Value = 10145.714
$ws.Range("M77").Value = -5777.714

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 29750
$ws.Range("J59").Value = 29750
$ws.Range("L59").Value = 29750
$ws.Range("N59").Value = -31444

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1390679.5
$ws.Range("I31").Value = 1725509.9
$ws.Range("J31").Value = 3524.9285
$ws.Range("K31").Value = 1725509.9
$ws.Range("L31").Value = 3524.9285
$ws.Range("M31").Value = -1725214.9
$ws.Range("N31").Value = -4114.9285

$ws.Range("H34").Value = 1390679.5
$ws.Range("I34").Value = 1725509.9
$ws.Range("J34").Value = 3524.9285
$ws.Range("K34").Value = 1725509.9
$ws.Range("L34").Value = 3524.9285
$ws.Range("M34").Value = -1725307.9
$ws.Range("N34").Value = -3928.9285

$ws.Range("H58").Value = 13517507
$ws.Range("I58").Value = 2149.087
$ws.Range("J58").Value = 35721308
$ws.Range("K58").Value = 2149.087
$ws.Range("L58").Value = 35721308
$ws.Range("M58").Value = -1946.087
$ws.Range("N58").Value = -35721714

$ws.Range("H99").Value = 2253.0527
$ws.Range("I99").Value = 1878
$ws.Range("J99").Value = 2525.818
$ws.Range("K99").Value = 1878
$ws.Range("L99").Value = 2525.818
$ws.Range("M99").Value = -380
$ws.Range("N99").Value = -5521.818

$ws.Range("H126").Value = 2253.0527
$ws.Range("I126").Value = 1878
$ws.Range("J126").Value = 2525.818
$ws.Range("K126").Value = 5634
$ws.Range("L126").Value = 7577.454000000001
$ws.Range("M126").Value = -3164
$ws.Range("N126").Value = -12517.454

$ws.Range("H132").Value = 2223.95
$ws.Range("I132").Value = 1922.72
$ws.Range("J132").Value = 2726
$ws.Range("K132").Value = 5768.16
$ws.Range("L132").Value = 8178
$ws.Range("M132").Value = -3238.16
$ws.Range("N132").Value = -13238

$ws.Range("H134").Value = 2683.8125
$ws.Range("I134").Value = 1099.8889
$ws.Range("J134").Value = 4720.2856
$ws.Range("K134").Value = 3299.6667
$ws.Range("L134").Value = 14160.8568
$ws.Range("M134").Value = -764.6666999999998
$ws.Range("N134").Value = -19230.8568

$ws.Range("H136").Value = 13517507
$ws.Range("I136").Value = 2149.087
$ws.Range("J136").Value = 35721308
$ws.Range("K136").Value = 6447.261
$ws.Range("L136").Value = 107163924
$ws.Range("M136").Value = -3897.261
$ws.Range("N136").Value = -107169024

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 1500
$ws.Range("J36").Value = 1500
$ws.Range("L36").Value = 4500
$ws.Range("N36").Value = -4838

$ws.Range("H37").Value = 49120
$ws.Range("J37").Value = 49120
$ws.Range("L37").Value = 147360
$ws.Range("N37").Value = -147584

$ws.Range("H68").Value = 3566.8333
$ws.Range("J68").Value = 4555.778
$ws.Range("L68").Value = 13667.334
$ws.Range("N68").Value = -15289.334

$ws.Range("H71").Value = 3566.8333
$ws.Range("J71").Value = 4555.778
$ws.Range("L71").Value = 41002.002
$ws.Range("N71").Value = -49114.002

$ws.Range("H131").Value = 1520.25
$ws.Range("J131").Value = 1260.2963
$ws.Range("L131").Value = 3780.8889
$ws.Range("N131").Value = -13860.8889

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4184.7095
$ws.Range("I132").Value = 4429.3335
$ws.Range("J132").Value = 3846
$ws.Range("K132").Value = 13288.0005
$ws.Range("L132").Value = 11538
$ws.Range("M132").Value = -10758.0005
$ws.Range("N132").Value = -16598

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2204.3333
$ws.Range("I132").Value = 1330.0435
$ws.Range("J132").Value = 3461.125
$ws.Range("K132").Value = 3990.1305
$ws.Range("L132").Value = 10383.375
$ws.Range("M132").Value = -1460.1305
$ws.Range("N132").Value = -15443.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 669152.75
$ws.Range("I122").Value = 1252060.8
$ws.Range("J122").Value = 2972.1428
$ws.Range("K122").Value = 3756182.4
$ws.Range("L122").Value = 8916.428400000001
$ws.Range("M122").Value = -3753732.4
$ws.Range("N122").Value = -13816.4284

$ws.Range("H132").Value = 317985.44
$ws.Range("I132").Value = 834807.7
$ws.Range("J132").Value = 7892.1
$ws.Range("K132").Value = 2504423.1
$ws.Range("L132").Value = 23676.3
$ws.Range("M132").Value = -2501893.1
$ws.Range("N132").Value = -28736.3
